# Append two new benchmark rows (32-33) to the log_evaluations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = "mistral:7b-instruct-v0.3-q5_K_M"
$ws.Range("B32").Value = "llama3:70b"
$ws.Range("C32").Value = 10
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = 404.9
$ws.Range("F32").Value = 20.81
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_match_4.txt"
$ws.Range("I32").Value = 20.81
$ws.Range("J32").Value = 1
$ws.Range("K32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_correct_4.txt"
$ws.Range("L32").Value = 20.81
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_executable_4.txt"
$ws.Range("O32").Value = 0.8571428571428571
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 22.1
$ws.Range("R32").Value = 0.5
$ws.Range("S32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_match_4.txt"
$ws.Range("T32").Value = 22.1
$ws.Range("U32").Value = 0.5
$ws.Range("V32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_correct_4.txt"
$ws.Range("W32").Value = 22.1
$ws.Range("X32").Value = 0.5
$ws.Range("Y32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_executable_4.txt"
$ws.Range("Z32").Value = 0
$ws.Range("AA32").Value = 334.56
$ws.Range("AB32").Value = 27.42
$ws.Range("AC32").Value = 0
$ws.Range("AD32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_match_4.txt"
$ws.Range("AE32").Value = 27.42
$ws.Range("AF32").Value = 0.5
$ws.Range("AG32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_correct_4.txt"
$ws.Range("AH32").Value = 27.42
$ws.Range("AI32").Value = 0.5
$ws.Range("AJ32").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_executable_4.txt"
$ws.Range("AK32").Value = 0
$ws.Range("AL32").Value = 2
$ws.Range("AM32").Value = 2
$ws.Range("AN32").Value = "text"
$ws.Range("AO32").Value = 140
$ws.Range("AP32").Value = 0.1
$ws.Range("AQ32").Value = 300
$ws.Range("AR32").Value = 0.9
$ws.Range("AS32").Value = 5
$ws.Range("AT32").Value = 1
$ws.Range("AU32").Value = 1.2
$ws.Range("AV32").Value = 1
$ws.Range("AW32").Value = 1024
$ws.Range("AX32").Value = ""

# Row 33
$ws.Range("A33").Value = "mistral:7b-instruct-v0.3-q5_K_M"
$ws.Range("B33").Value = "llama3:70b"
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 200
$ws.Range("E33").Value = 404.32
$ws.Range("F33").Value = 20.95
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_match_4.txt"
$ws.Range("I33").Value = 20.95
$ws.Range("J33").Value = 1
$ws.Range("K33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_correct_4.txt"
$ws.Range("L33").Value = 20.95
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_executable_4.txt"
$ws.Range("O33").Value = 0.8571428571428571
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 21.66
$ws.Range("R33").Value = 0.5
$ws.Range("S33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_match_4.txt"
$ws.Range("T33").Value = 21.66
$ws.Range("U33").Value = 0.5
$ws.Range("V33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_correct_4.txt"
$ws.Range("W33").Value = 21.66
$ws.Range("X33").Value = 0.5
$ws.Range("Y33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_executable_4.txt"
$ws.Range("Z33").Value = 0
$ws.Range("AA33").Value = 339.11
$ws.Range("AB33").Value = 22.59
$ws.Range("AC33").Value = 0.5
$ws.Range("AD33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_match_4.txt"
$ws.Range("AE33").Value = 22.59
$ws.Range("AF33").Value = 1
$ws.Range("AG33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_correct_4.txt"
$ws.Range("AH33").Value = 22.59
$ws.Range("AI33").Value = 1
$ws.Range("AJ33").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_executable_4.txt"
$ws.Range("AK33").Value = 0.4285714285714285
$ws.Range("AL33").Value = 2
$ws.Range("AM33").Value = 2
$ws.Range("AN33").Value = "text"
$ws.Range("AO33").Value = 140
$ws.Range("AP33").Value = 0.1
$ws.Range("AQ33").Value = 300
$ws.Range("AR33").Value = 0.9
$ws.Range("AS33").Value = 5
$ws.Range("AT33").Value = 1
$ws.Range("AU33").Value = 1.2
$ws.Range("AV33").Value = 1
$ws.Range("AW33").Value = 1024
$ws.Range("AX33").Value = ""

